$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List_ID")

# Fix the pass/fail values on row 3 (swap)
$ws.Range("F3").Value = "fail"
$ws.Range("G3").Value = "PASS"

# Append new data rows 18-21
# NOTE: "09/01/2001" is ambiguous as a US-style m/d/yyyy date (month=09, day=01)
# so a plain .Value assignment gets auto-parsed into a date serial number.
# Force it to stay literal text by staging it through a Text-formatted scratch
# cell and pasting only the value back in, then removing the scratch cell.
$ws.Range("ZZ1").NumberFormat = "@"
$ws.Range("ZZ1").Value = "09/01/2001"
$ws.Range("ZZ1").Copy()
$ws.Range("C18").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

$ws.Range("D18").Value = "opqsleciiggdhik@gmail.com"
$ws.Range("E18").Value = "wdmujKSMZF5"
$ws.Range("F18").Value = "pass"

$ws.Range("C19").Value = "17/07/1992"
$ws.Range("D19").Value = "bbchdfmnfrcpkgm@gmail.com"
$ws.Range("E19").Value = "dsmtzHUSSJ5"
$ws.Range("F19").Value = "pass"

$ws.Range("C20").Value = "28/09/1974"
$ws.Range("D20").Value = "oqvbovkticuqkqb@gmail.com"
$ws.Range("E20").Value = "tepfsZSOFA5"
$ws.Range("F20").Value = "pass"

$ws.Range("C21").Value = "24/06/1992"
$ws.Range("D21").Value = "ccjatqnecgvwuey@gmail.com"
$ws.Range("E21").Value = "xdtpcXQAEI5"
$ws.Range("F21").Value = "pass"

# Update the selection to reflect the saved cursor position
$ws.Range("F4").Select()
